$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in rows 2-5 (A and B columns) to new data
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 46

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 21

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 20

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 12

# Delete row 6 entirely (shifts nothing up, just removes last row's data)
$ws.Range("A6:B6").Delete()
